$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "studies": rename/restructure header row + row 2 data
# ----------------------------------------------------------------------
$studies = $wb.Worksheets.Item("studies")

$studies.Range("A1").Value = "study_id"
$studies.Range("B1").Value = "study_label"
$studies.Range("C1").Value = "description"
$studies.Range("D1").Value = "access_level"
$studies.Range("E1").Value = "contributors"
$studies.Range("F1").Value = "reference"
$studies.Range("G1").Value = "reference_year"

# ----------------------------------------------------------------------
# Sheet "surveys": rewrite the header/data columns for the new location_* /
# time_* fields. Cells are written directly (no column insert/shift) because
# the new "collection_*" columns keep living in I:K while the worksheet's
# <cols style="2"> range is left pointing at the original 8:10 span.
# ----------------------------------------------------------------------
$surveys = $wb.Worksheets.Item("surveys")

$surveys.Range("A1").Value = "study_id"

# old G (spatial_notes) becomes location_method, its data value moves to the
# brand-new "location_notes" column (H)
$surveys.Range("G2").ClearContents()
$surveys.Range("G1").Value = "location_method"
$surveys.Range("H1").Style = "Normal"
$surveys.Range("H1").Value = "location_notes"
$surveys.Range("H2").Value = "example data"
$surveys.Range("H2").Style = "Normal"

# collection_start/end/day shift from H:J to I:K, keeping the "@" text format
$surveys.Range("I1").NumberFormat = "@"
$surveys.Range("I1").Value = "collection_start"
$surveys.Range("J1").NumberFormat = "@"
$surveys.Range("J1").Value = "collection_end"
$surveys.Range("K1").NumberFormat = "@"
$surveys.Range("K1").Value = "collection_day"
$surveys.Range("J2").Style = "Normal"
$surveys.Range("J2").ClearContents()
$surveys.Range("K2").NumberFormat = "@"
$surveys.Range("K2").Value = "2020-01-01"

# a new "time_method" column, and the trailing time_notes column moves from K to M
$surveys.Range("L1").NumberFormat = "@"
$surveys.Range("L1").Value = "time_method"
$surveys.Range("L2").NumberFormat = "@"
$surveys.Range("M1").Style = "Normal"
$surveys.Range("M1").Value = "time_notes"
$surveys.Range("M2").Value = "example data"

# Row 2 of studies: data moves around - C2 (study_type "other") disappears, D2
# becomes "public", F2 keeps the hyperlink text (now under "reference"), G2 stays blank.
$studies.Range("C2").ClearContents()
$studies.Range("D2").Value = "public"

# ----------------------------------------------------------------------
# Sheet "counts": rename key columns (values are unchanged)
# ----------------------------------------------------------------------
$counts = $wb.Worksheets.Item("counts")
$counts.Range("A1").Value = "study_id"
$counts.Range("B1").Value = "survey_id"

# ----------------------------------------------------------------------
# Selections / active sheet
# ----------------------------------------------------------------------
$studies.Range("D8").Select()
$surveys.Range("C6").Select()
$counts.Range("B2").Select()

$studies.Activate()
$studies.Range("D8").Select()
